# Build site at 2023-04-12 14:53:07 UTC
# Updates the LOM3032 "Produtos Ceramicos" syllabus sheet:
#  - Objetivos/Objectives now carries real course-objective text, and the
#    teacher-name row ("1922320 - Sebastiao Ribeiro") moves down to sit
#    under "Docentes responsaveis:" in a brand-new row.
#  - Programa resumido / Programa get real syllabus text (replacing the
#    placeholder "Semestral" / date values).
#  - Criterio / Norma de recuperacao / Bibliografia get their real text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new row at 13 - this pushes the old rows 13..21 down to 14..22,
#    carrying their styles/heights along for the ride.
$ws.Rows.Item(13).Insert()

# The insert leaves a style-only leftover in A13 (copied down from A12);
# the target layout has no A13 cell at all, so drop it.
$ws.Range("A13").Clear()

# Give the new B13/C13 cells the same look as the rest of column B/C
# (vertical-top, wrap-text) by lifting the format from row 10, then fill
# in the value that used to live in the old row-13/row-18 slot.
$ws.Range("B10:C10").Copy()
$ws.Range("B13:C13").PasteSpecial(-4122)
$ws.Range("B13").Value = "1922320 - Sebastiao Ribeiro"
$ws.Range("C13").Value = "1922320 - Sebastiao Ribeiro"

# 2) Objetivos: row now holds the real objective text instead of the
#    teacher name (which moved to row 13 above).
$ws.Range("B10").Value = "Informar os alunos dos produtos cerâmicos, de suas propriedades e suas aplicações"
$ws.Range("C10").Value = "Informar os alunos dos produtos cerâmicos, de suas propriedades e suas aplicações"

# 3) Programa resumido: (row 14, was row 13) - real short syllabus text
#    instead of the placeholder "Semestral".
$ws.Range("B14").Value = "1.Cerâmica vermelha2.Cerâmica branca3.Refratários4.Cerâmica eletro-eletrônica5.Cerâmica avançada estrutural"
$ws.Range("C14").Value = "1.Cerâmica vermelha2.Cerâmica branca3.Refratários4.Cerâmica eletro-eletrônica5.Cerâmica avançada estrutural"

# 4) Programa: (row 16, was row 15) - real syllabus text instead of the
#    placeholder date "01/01/1996".
$ws.Range("B16").Value = "01 - Cerâmica vermelha02 - Cerâmica vermelha03 - Cerâmica vermelha04 - Cerâmica vermelha05 - Cerâmica branca06 - Cerâmica branca07 - Refratários08 - Refratários09 - Refratários10 - Cerâmica eletro-eletrônica11 - Cerâmica eletro-eletrônica12 - Cerâmica eletro-eletrônica13 - Cerâmica eletro-eletrônica14 - Cerâmica avançada estrutural15 - Cerâmica avançada estrutural"
$ws.Range("C16").Value = "01 - Cerâmica vermelha02 - Cerâmica vermelha03 - Cerâmica vermelha04 - Cerâmica vermelha05 - Cerâmica branca06 - Cerâmica branca07 - Refratários08 - Refratários09 - Refratários10 - Cerâmica eletro-eletrônica11 - Cerâmica eletro-eletrônica12 - Cerâmica eletro-eletrônica13 - Cerâmica eletro-eletrônica14 - Cerâmica avançada estrutural15 - Cerâmica avançada estrutural"

# 5) Metodo: (row 19, was row 18) - real grading-method text instead of
#    the stray teacher-name value.
$ws.Range("B19").Value = "Duas provas escritas (P1 e P2), valendo de 0 (zero) a 10 (dez)"
$ws.Range("C19").Value = "Duas provas escritas (P1 e P2), valendo de 0 (zero) a 10 (dez)"

# 6) Criterio: (row 20, was row 19) - real pass/fail criteria text.
$ws.Range("B20").Value = "Média Parcial (MP): (P1 + P2)/2Média Parcial igual ou superior a 5: aprovação diretaMédia Parcial entre 3 e 5: recuperaçãoMédia Parcial inferior a 5: reprovação direta"
$ws.Range("C20").Value = "Média Parcial (MP): (P1 + P2)/2Média Parcial igual ou superior a 5: aprovação diretaMédia Parcial entre 3 e 5: recuperaçãoMédia Parcial inferior a 5: reprovação direta"

# 7) Norma de recuperacao: (row 21, was row 20) - real makeup-exam text.
$ws.Range("B21").Value = "Para a recuperação será realizada uma prova (PR) abrangendo toda a matéria lecionada no semestre, valendo de 0 (zero) a 10 (dez)Média Final:(MP + PR)/2Média Final igual ou superior a 5 (cinco): aprovadoMédia Final inferior a 5: reprovado"
$ws.Range("C21").Value = "Para a recuperação será realizada uma prova (PR) abrangendo toda a matéria lecionada no semestre, valendo de 0 (zero) a 10 (dez)Média Final:(MP + PR)/2Média Final igual ou superior a 5 (cinco): aprovadoMédia Final inferior a 5: reprovado"

# 8) Bibliografia: (row 22, was row 21) - real bibliography text instead
#    of the makeup-exam text that now lives one row up.
$ws.Range("B22").Value = "1.F. Singer & S. S. Singer, Cerâmica Industrial, V. 11, 19712.Salmang & Scholze, Keramik: Teil2 Keramische Werkstoffe, Springer Verlag, 19833.L. M. Levinson, Electronic Ceramics, Properties, Devices and Applications4.M. J. Hoffmann, Silicon Nitride"
$ws.Range("C22").Value = "1.F. Singer & S. S. Singer, Cerâmica Industrial, V. 11, 19712.Salmang & Scholze, Keramik: Teil2 Keramische Werkstoffe, Springer Verlag, 19833.L. M. Levinson, Electronic Ceramics, Properties, Devices and Applications4.M. J. Hoffmann, Silicon Nitride"
